$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 1526
$ws.Range("J3").Value = 1596
$ws.Range("I4").Value = 1756
$ws.Range("J4").Value = 362
$ws.Range("J6").Value = 2079
$ws.Range("I7").Value = 26200
$ws.Range("J7").Value = 5673

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J3").Value = 14
$ws.Range("J7").Value = 65

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 21

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 55
$ws.Range("J3").Value = 66
$ws.Range("J6").Value = 60
$ws.Range("J7").Value = 194

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J2").Value = 21
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 72

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 40
$ws.Range("J3").Value = 84
$ws.Range("J6").Value = 62
$ws.Range("J7").Value = 201

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 39
$ws.Range("J6").Value = 63
$ws.Range("J7").Value = 144

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J5").Value = 15
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 160
$ws.Range("J8").Value = 348
$ws.Range("J11").Value = 70
$ws.Range("J15").Value = 72
$ws.Range("J16").Value = 19
$ws.Range("J18").Value = 69
$ws.Range("J20").Value = 118
$ws.Range("J23").Value = 51
$ws.Range("J24").Value = 21
$ws.Range("J25").Value = 32
$ws.Range("J27").Value = 33
$ws.Range("J29").Value = 317
$ws.Range("J30").Value = 21
$ws.Range("J33").Value = 236
$ws.Range("J36").Value = 87
$ws.Range("J37").Value = 194
$ws.Range("J41").Value = 36
$ws.Range("J42").Value = 222
$ws.Range("J53").Value = 53
$ws.Range("J57").Value = 26
$ws.Range("J60").Value = 32
$ws.Range("I63").Value = 193
$ws.Range("J63").Value = 24
$ws.Range("J65").Value = 144
$ws.Range("J67").Value = 201
$ws.Range("J71").Value = 27
$ws.Range("J73").Value = 55
$ws.Range("J79").Value = 175
$ws.Range("J83").Value = 140
$ws.Range("J85").Value = 257
$ws.Range("J88").Value = 58
$ws.Range("J89").Value = 65
$ws.Range("J92").Value = 16
$ws.Range("J94").Value = 43
$ws.Range("J95").Value = 83
$ws.Range("J97").Value = 37
$ws.Range("J99").Value = 72
$ws.Range("I101").Value = 26200
$ws.Range("J101").Value = 5673

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J6").Value = 42
$ws.Range("J7").Value = 140

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 83

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 60
$ws.Range("J3").Value = 68
$ws.Range("J6").Value = 91
$ws.Range("J7").Value = 236

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 93
$ws.Range("J3").Value = 118
$ws.Range("J6").Value = 80
$ws.Range("J7").Value = 317

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 62
$ws.Range("J3").Value = 102
$ws.Range("J6").Value = 70
$ws.Range("J7").Value = 257

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J2").Value = 18
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 56

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J6").Value = 16
$ws.Range("J7").Value = 36

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 47
$ws.Range("J7").Value = 222

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 21

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 47
$ws.Range("J3").Value = 62
$ws.Range("J6").Value = 50
$ws.Range("J7").Value = 175

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J2").Value = 31
$ws.Range("J3").Value = 38
$ws.Range("J7").Value = 118

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J2").Value = 19
$ws.Range("J7").Value = 69

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J2").Value = 30
$ws.Range("J7").Value = 87

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J6").Value = 26
$ws.Range("J7").Value = 43

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J2").Value = 14
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 32

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 72

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J6").Value = 24
$ws.Range("J7").Value = 70

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 55

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J6").Value = 25
$ws.Range("J7").Value = 37

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 16

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J6").Value = 30
$ws.Range("J7").Value = 58

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J3").Value = 114
$ws.Range("J4").Value = 18
$ws.Range("J6").Value = 95
$ws.Range("J7").Value = 348

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 15

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 26

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 32

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J3").Value = 13
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("J2").Value = 6
$ws.Range("J3").Value = 5
$ws.Range("J7").Value = 27

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J3").Value = 52
$ws.Range("J6").Value = 54
$ws.Range("J7").Value = 160

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 19
